# Estadisticos Matutinos 15 Oct
# Updates statistics rows (rows 3-5) on "Estadisticos 1P" and "Estadisticos Final"
# sheets with Reprobados/Aprobados/Por_Apro/Promedio figures, updates the
# "Aprobados" column (E) on "Estadisticos 2P", and appends the list of
# rescatable (recoverable) students on the "Rescatables" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Estadisticos 1P" (sheet1) and "Estadisticos Final" (sheet3): rows 3-5,
# columns D (Reprobados), F (Aprobados), G (Por_Apro), H (Promedio)
# ---------------------------------------------------------------------------
$statsData = @(
    @{ Row = 3; D = 13; F = 31; G = 70.45; H = 7.1 },
    @{ Row = 4; D = 13; F = 30; G = 69.77; H = 7 },
    @{ Row = 5; D = 6;  F = 12; G = 66.67; H = 6.9 }
)

foreach ($sheetName in @("Estadisticos 1P", "Estadisticos Final")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $statsData) {
        $r = $entry.Row
        $ws.Cells.Item($r, 4).Value = $entry.D
        $ws.Cells.Item($r, 6).Value = $entry.F
        $ws.Cells.Item($r, 7).Value = $entry.G
        $ws.Cells.Item($r, 8).Value = $entry.H
    }
}

# ---------------------------------------------------------------------------
# "Estadisticos 2P" (sheet2): rows 3-5, column E (Reprobados)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Cells.Item(3, 5).Value = 31
$ws2.Cells.Item(4, 5).Value = 30
$ws2.Cells.Item(5, 5).Value = 12

# ---------------------------------------------------------------------------
# "Rescatables" (sheet4): append student rows 2-9
# Columns: A=NC, B=Paterno, C=Materno, D=Nombres, E=Nombre_Largo, F=Grupo, G=Reprobadas
# ---------------------------------------------------------------------------
$wsR = $wb.Worksheets.Item("Rescatables")

$rescatables = @(
    @{ NC = 21330051920005; Paterno = "CASTILLO";   Materno = "GARCIA";    Nombres = "KEVIN ISAAC";    Grupo = "1AV"; Reprobadas = 6 },
    @{ NC = 21330051920078; Paterno = "GARCIA";      Materno = "GONZALEZ";  Nombres = "MIROSLAVA";      Grupo = "1CV"; Reprobadas = 6 },
    @{ NC = 21330051920092; Paterno = "OLIVARES";    Materno = "HIPOLITO";  Nombres = "JOSE JULIAN";    Grupo = "1CV"; Reprobadas = 6 },
    @{ NC = 21330051920098; Paterno = "RAMIREZ";     Materno = "MARTINEZ";  Nombres = "ALDO GEOVANNI";  Grupo = "1CV"; Reprobadas = 6 },
    @{ NC = 21330051920105; Paterno = "MENA";        Materno = "ANGELES";   Nombres = "SOLANO";         Grupo = "1CV"; Reprobadas = 6 },
    @{ NC = 21330051920135; Paterno = "ALVAREZ";     Materno = "RIVERA";    Nombres = "PEDRO ANGEL";    Grupo = "1EV"; Reprobadas = 6 },
    @{ NC = 21330051920136; Paterno = "ARELLANO";    Materno = "JUAREZ";    Nombres = "DAVID OSWALDO";  Grupo = "1EV"; Reprobadas = 6 },
    @{ NC = 21330051920038; Paterno = "DOMINGUEZ";   Materno = "TORRES";    Nombres = "ZAYRA";          Grupo = "1BV"; Reprobadas = 6 }
)

$row = 2
foreach ($s in $rescatables) {
    $wsR.Cells.Item($row, 1).Value = $s.NC
    $wsR.Cells.Item($row, 2).Value = $s.Paterno
    $wsR.Cells.Item($row, 3).Value = $s.Materno
    $wsR.Cells.Item($row, 4).Value = $s.Nombres
    $wsR.Cells.Item($row, 5).Value = "INGLÉS I"
    $wsR.Cells.Item($row, 6).Value = $s.Grupo
    $wsR.Cells.Item($row, 7).Value = $s.Reprobadas
    $row++
}
